$d = $word.ActiveDocument

# Locate the paragraph ending in "...рябиновую" (Усачев Илья bullet item)
# so the new bullet item gets inserted right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*рябиновую*") {
        $target = $p
        break
    }
}

$r = $target.Range
$r.InsertParagraphAfter()

# The newly created paragraph inherits the list/border/shading formatting
# from the paragraph it followed; fill in its text.
$newPara = $target.Next()
$newPara.Range.InsertAfter("Ольга Маймасова – логист. Прямой начальник операторов ММЗ, в случае если он долго не выходит на связь или пропал вовсе – пишем звоним ей!")
